# edit.ps1 -- apply the "before merging Hisan's project" change set
#
# Summary of changes:
#  1. Three empty table-cell paragraphs (the blank "headline" cells in the
#     first table) lose their stray <w:rFonts w:hint="eastAsia"/> paragraph
#     mark formatting and become plain empty paragraphs.
#  2. The stray _GoBack bookmark sitting inside the "Answer_pred" header
#     cell is removed.
#  3. A new bulleted ("Pipeline" style numbered list, numId 7) paragraph is
#     appended at the very end of the document body with the WordSim note,
#     and the _GoBack bookmark is re-created there (i.e. the bookmark just
#     "moved" to the newest edit, exactly like Word does automatically).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Clean up the three blank cells in the first ("Story_id / headline /
#    date / Story type") table -- column 2, rows 2-4.
# ---------------------------------------------------------------------
$tbl1 = $d.Tables(1)
for ($r = 2; $r -le 4; $r++) {
    $cell = $tbl1.Cell($r, 2)
    $p = $cell.Range.Paragraphs(1)
    $p.Range.InsertXML("<w:p $wNs/>") | Out-Null
}

# ---------------------------------------------------------------------
# 2) Strip the _GoBack bookmark out of the "Answer_pred" header cell
#    (2nd table, row 1, last column) while keeping its text/runs intact.
# ---------------------------------------------------------------------
$tbl2 = $d.Tables(2)
$answerPredCell = $tbl2.Cell(1, $tbl2.Columns.Count)
$answerPredPara = $answerPredCell.Range.Paragraphs(1)
$cleanCellXml = @"
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:color w:val="FF0000"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:color w:val="FF0000"/>
    </w:rPr>
    <w:t>Ans</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:color w:val="FF0000"/>
    </w:rPr>
    <w:t>wer</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:color w:val="FF0000"/>
    </w:rPr>
    <w:t>_pred</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
$answerPredPara.Range.InsertXML($cleanCellXml) | Out-Null

# ---------------------------------------------------------------------
# 3) Append the new "WordSim(...)" bullet paragraph at the end of the
#    document (same list as the "Pipeline" bullet, numId 7) and drop the
#    _GoBack bookmark there.
# ---------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()

$newParaIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newParaIndex)

$newParaXml = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="7"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>WordSim</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>work, works)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri" w:hint="eastAsia"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>=</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri" w:hint="eastAsia"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>0.83!</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Microsoft YaHei" w:cs="Calibri"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> We need to compare the lemmas!!!</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@
$newPara.Range.InsertXML($newParaXml) | Out-Null
